# Controle da presenca (dia 19/09/2022) e notas do forum para a semana
# 11/09/2022 a 17/09/2022: recalcula a coluna J (nota_view), reduzindo de
# 5 para 4 as notas maximas atribuidas anteriormente.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 10)  # column J = 10
    if ($cell.Value2 -eq 5) {
        $cell.Value2 = 4
    }
}
